$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update data row 3 (Excel row 3) values per the commit:
# "se agregó al TC: Emision_Motor 2 recording para que seleccione el grupo del PAS"
$ws.Range("A3").Value = "i-preproducciongestion.segurossura.com.ar"
$ws.Range("E3").Value = 7105947778
$ws.Range("I3").Value = "Cupón"
$ws.Range("P3").Value = "CADETE EN MOTO"
$ws.Range("T3").Value = 21004018

# Update the active window's view to reflect the new selection / scroll position
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = $ws.Range("K1").Column
$ws.Range("T9").Select()
